$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (GitHub Actions scheduled update)

$ws.Range('D2').Value = '64.933.82'
$ws.Range('E2').Value = '  +5.72%  '

$ws.Range('D3').Value = '2.975.72'
$ws.Range('E3').Value = '  +2.96%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.01'
$ws.Range('E5').Value = '  +2.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.18'
$ws.Range('E6').Value = '  +7.24%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  +1.17%  '

$ws.Range('D9').Value = '2.972.54'
$ws.Range('E9').Value = '  +2.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.96'
$ws.Range('E10').Value = '  +3.57%  '

$ws.Range('E11').Value = '  +2.08%  '

$ws.Range('E12').Value = '  +2.89%  '

$ws.Range('E13').Value = '  +1.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.83'
$ws.Range('E14').Value = '  +5.98%  '

$ws.Range('E15').Value = '  +0.70%  '

$ws.Range('D16').Value = '64.869.61'
$ws.Range('E16').Value = '  +5.52%  '

$ws.Range('D17').Value = '3.468.21'
$ws.Range('E17').Value = '  +2.93%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.87'
$ws.Range('E18').Value = '  +3.56%  '

$ws.Range('D19').Value = '2.979.54'
$ws.Range('E19').Value = '  +3.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '447.56'
$ws.Range('E20').Value = '  +3.52%  '

$ws.Range('E21').Value = '  +3.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.677'
$ws.Range('E22').Value = '  +3.18%  '

$ws.Range('E23').Value = '  +5.00%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.80'
$ws.Range('E24').Value = '  +1.78%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.23'
$ws.Range('E25').Value = '  +3.70%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.55'
$ws.Range('E26').Value = '  +5.53%  '

$ws.Range('E27').Value = '  +7.64%  '

$ws.Range('E28').Value = '  -0.06%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.78'
$ws.Range('E29').Value = '  +11.15%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.35'
$ws.Range('E30').Value = '  +14.27%  '

$ws.Range('E31').Value = '  +2.66%  '

$ws.Range('E32').Value = '  -1.33%  '

$ws.Range('E33').Value = '  +3.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.51'
$ws.Range('E34').Value = '  +3.44%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.17%  '

$ws.Range('E36').Value = '  +2.49%  '

$ws.Range('E37').Value = '  +4.01%  '

$ws.Range('E38').Value = '  +7.36%  '

$ws.Range('E39').Value = '  +0.29%  '

$ws.Range('E40').Value = '  +1.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '43.63'
$ws.Range('E41').Value = '  +10.05%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  +5.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.295'
$ws.Range('E43').Value = '  +10.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.38'
$ws.Range('E44').Value = '  +1.67%  '

$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '378.57'
$ws.Range('E45').Value = '  +11.51%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.760.06'
$ws.Range('E46').Value = '  +2.83%  '

$ws.Range('E47').Value = '  +3.72%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.78'
$ws.Range('E48').Value = '  -0.16%  '

$ws.Range('E49').Value = '  -0.04%  '

$ws.Range('E50').Value = '  +2.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.83'
$ws.Range('E51').Value = '  +6.10%  '
